$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = '29.020.42'
$ws.Range("E2").Value = '  -0.77%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = '1.828.70'
$ws.Range("E3").Value = '  -0.78%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.28'
$ws.Range("D5").ClearFormats()

# Row 6: 'XRP' -> 'XRP'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6532'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.15%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = '  +0.05%  '

# Row 8: 'OKB' -> 'OKB'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.43'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.58%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2935'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.47%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07335'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.24%  '

# Row 11: 'Solana' -> 'Solana'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.91'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.35%  '

# Row 12: 'TRON' -> 'TRON'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07668'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.63%  '

# Row 13: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D13").Value = '1.834.52'
$ws.Range("E13").Value = '  -0.90%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.978'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.46%  '

# Row 15: 'Polygon' -> 'Polygon'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6664'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.58%  '

# Row 16: 'Litecoin' -> 'Litecoin'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.84'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.90%  '

# Row 17: 'Uniswap' -> 'Uniswap'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.097'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.40%  '

# Row 18: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008655'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.11%  '

# Row 19: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D19").Value = '29.007.01'
$ws.Range("E19").Value = '  -0.95%  '

# Row 20: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D20").Value = '2.084.80'
$ws.Range("E20").Value = '  -0.93%  '

# Row 21: 'Avalanche' -> 'Avalanche'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.63%  '

# Row 22: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.92'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.02%  '

# Row 23: 'Dai' -> 'Dai'
$ws.Range("E23").Value = '  -0.10%  '

# Row 24: 'Chainlink' -> 'Chainlink'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.100'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.23%  '

# Row 25: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.001'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.08%  '

# Row 26: 'Monero' -> 'Monero'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.94'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.73%  '

# Row 27: 'Cosmos' -> 'Cosmos'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.499'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.25%  '

# Row 28: 'Stellar' -> 'Stellar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1377'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.81%  '

# Row 29: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.86'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.74%  '

# Row 30: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.506'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.22%  '

# Row 31: 'Filecoin' -> 'Filecoin'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.106'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.68%  '

# Row 32: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.012'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.40%  '

# Row 33: 'Toncoin' -> 'Toncoin'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.200'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.81%  '

# Row 34: 'Hedera' -> 'Hedera'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05351'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.79%  '

# Row 35: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7429'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.45%  '

# Row 36: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.839'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.74%  '

# Row 37: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.154'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.68%  '

# Row 38: 'HuobiToken' -> 'HuobiToken'
$ws.Range("E38").Value = '  -1.22%  '

# Row 39: 'Maker' -> 'Maker'
$ws.Range("D39").Value = '1.299.31'

# Row 40: 'VeChain' -> 'VeChain'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01786'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.07%  '

# Row 41: 'MXToken' -> 'MXToken'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.749'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.97%  '

# Row 42: 'FraxShare' -> 'FraxShare'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.367'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +7.30%  '

# Row 43: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8979'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.33%  '

# Row 44: 'PaxDollar' -> 'PaxDollar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.26%  '

# Row 45: 'Quant' -> 'Quant'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.15'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.33%  '

# Row 46: 'RocketPoolETH' -> 'RocketPoolETH'
$ws.Range("D46").Value = '1.983.81'
$ws.Range("E46").Value = '  -0.78%  '

# Row 47: 'Mantle' -> 'Mantle'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5141'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.32%  '

# Row 48: 'Aave' -> 'Aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.11'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.55%  '

# Row 49: 'BabyDogeCoin' -> 'RenderToken'
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.739'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.17%  '

# Row 50: 'XinFinNetwork' -> 'XinFinNetwork'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07535'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -7.35%  '

# Row 51: 'RenderToken' -> 'BabyDogeCoin'
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000119'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.54%  '
